$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear removed cells
$ws.Range("F3").ClearContents()

# Cells that must remain text even though they look numeric
$c = $ws.Range("B2")
$c.NumberFormat = "@"
$c.Value = "5"
$c.Style = "Normal"

# Updated values
$ws.Range("D2").Value = 0.0354
$ws.Range("E2").Value = 0.1495
$ws.Range("F2").Value = -0.0207
$ws.Range("G2").Value = 0.3224429426102371
$ws.Range("H2").Value = 0.3224429426102371
$ws.Range("I2").Value = 0.2464181254154664
$ws.Range("J2").Value = 0.1859234386920439
$ws.Range("K2").Value = 1552.218
$ws.Range("L2").Value = 0.1719718590737868
$ws.Range("M2").Value = 797.8288
$ws.Range("N2").Value = 0.04423314427645549
$ws.Range("O2").Value = 0.5139927510182204
$ws.Range("P2").Value = 699.8008
$ws.Range("Q2").Value = 0.03879828573646247
$ws.Range("R2").Value = 0.4508392506722638
$ws.Range("S2").Value = 98.02799999999999
$ws.Range("T2").Value = 0.1228684650140481
$ws.Range("U2").Value = 641.337
$ws.Range("V2").Value = 0.03555694160304709
$ws.Range("W2").Value = 0.2876064333017975
$ws.Range("X2").Value = 0.05997947110027935
$ws.Range("Y2").Value = 0.2276269622015182
$ws.Range("Z2").Value = 1.924137069504708
$ws.Range("AA2").Value = 0.2454347283306264
$ws.Range("AB2").Value = 0.05978493928437828
$ws.Range("AC2").Value = 0.1847925174975664
$ws.Range("AD2").Value = 444.652
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 444.652
$ws.Range("AG2").Value = -196.685
$ws.Range("AH2").Value = 0.02405923485213796
$ws.Range("AI2").Value = 0.09286893435857335
$ws.Range("AJ2").Value = -0.01102481107991131
$ws.Range("AK2").Value = -0.04743266495683829
$ws.Range("AL2").Value = 32.919
$ws.Range("AM2").Value = 32.064
$ws.Range("AN2").Value = 0.200467070619635
$ws.Range("AO2").Value = 67.56493210607856
$ws.Range("AP2").Value = -0.08867353747385126
$ws.Range("AQ2").Value = 69.36657934131736
$ws.Range("B3").Value = "Wiz Soluções e Corretagem de Seguros S.A. (BOVESPA:WIZS3)"
$ws.Range("D3").Value = 0.159
$ws.Range("E3").Value = 0.19
$ws.Range("G3").Value = 0.7400749063670412
$ws.Range("H3").Value = 0.7400749063670412
$ws.Range("I3").Value = 0.501123595505618
$ws.Range("J3").Value = 0.3072174908752594
$ws.Range("K3").Value = 38
$ws.Range("L3").Value = 0.2846441947565543
$ws.Range("M3").Value = 26.5
$ws.Range("N3").Value = 0.1070274636510501
$ws.Range("O3").Value = 0.6973684210526315
$ws.Range("P3").Value = 26.5
$ws.Range("Q3").Value = 0.1070274636510501
$ws.Range("R3").Value = 0.6973684210526315
$ws.Range("U3").Value = 31.4
$ws.Range("V3").Value = 0.1268174474959612
$ws.Range("W3").Value = 0.6440677966101694
$ws.Range("X3").Value = 0.0599683918499534
$ws.Range("Y3").Value = 0.584099404760216
$ws.Range("Z3").Value = 40.45454545454549
$ws.Range("AA3").Value = 12.4283439490446
$ws.Range("AB3").Value = 0.05978493928437828
$ws.Range("AC3").Value = 12.36855900976022
$ws.Range("AD3").Value = 1.8
$ws.Range("AF3").Value = 1.8
$ws.Range("AG3").Value = -29.6
$ws.Range("AH3").Value = 0.007217321571772253
$ws.Range("AI3").Value = 0.02706766917293233
$ws.Range("AJ3").Value = -0.1357798165137614
$ws.Range("AK3").Value = -0.843304843304843
$ws.Range("AL3").Value = 0.201
$ws.Range("AM3").Value = -0.349
$ws.Range("AN3").Value = 0.0237467018469657
$ws.Range("AO3").Value = 332.8358208955224
$ws.Range("AP3").Value = -0.3905013192612137
$ws.Range("AQ3").Value = -191.6905444126074
$ws.Range("B4").Value = "BB Seguridade Participações S.A. (BOVESPA:BBSE3)"
$ws.Range("D4").Value = -0.00442
$ws.Range("E4").Value = -0.0111
$ws.Range("F4").Value = -0.101
$ws.Range("G4").Value = 1.335056440552801
$ws.Range("H4").Value = 1.335056440552801
$ws.Range("I4").Value = 0.9368076801350353
$ws.Range("J4").Value = 0.7438185694052909
$ws.Range("K4").Value = 729.6
$ws.Range("L4").Value = 0.7697014453001372
$ws.Range("M4").Value = 509.5
$ws.Range("N4").Value = 0.04472476057549662
$ws.Range("O4").Value = 0.698327850877193
$ws.Range("P4").Value = 509.5
$ws.Range("Q4").Value = 0.04472476057549662
$ws.Range("R4").Value = 0.698327850877193
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 325.3
$ws.Range("V4").Value = 0.02855537706616105
$ws.Range("W4").Value = 0.2876064333017975
$ws.Range("X4").Value = 0.05972644268231771
$ws.Range("Y4").Value = 0.2278799906194798
$ws.Range("Z4").Value = 0.886716557530402
$ws.Range("AA4").Value = 0.659556241290248
$ws.Range("AB4").Value = 0.05972644268231771
$ws.Range("AC4").Value = 0.5998297986079303
$ws.Range("AD4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -325.3
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.02939475539009271
$ws.Range("AK4").Value = -0.4034478481954608
$ws.Range("AL4").Value = 0.147
$ws.Range("AM4").Value = 0.147
$ws.Range("AN4").Value = 0
$ws.Range("AO4").Value = 6040.816326530613
$ws.Range("AP4").Value = -0.3663288288288288
$ws.Range("AQ4").Value = 6040.816326530613
$ws.Range("B5").Value = "Sul América S.A. (BOVESPA:SULA11)"
$ws.Range("D5").Value = 0.07780000000000001
$ws.Range("E5").Value = 0.302
$ws.Range("F5").Value = 0.00226
$ws.Range("G5").Value = 0.1959499702203693
$ws.Range("H5").Value = 0.1959499702203693
$ws.Range("I5").Value = 0.153685801988363
$ws.Range("J5").Value = 0.1106238538194302
$ws.Range("K5").Value = 490.9
$ws.Range("L5").Value = 0.1124524671278692
$ws.Range("M5").Value = 131.9
$ws.Range("N5").Value = 0.04011801204452826
$ws.Range("O5").Value = 0.2686901609289061
$ws.Range("P5").Value = 41
$ws.Range("Q5").Value = 0.01247034491149097
$ws.Range("R5").Value = 0.08352006518639235
$ws.Range("S5").Value = 90.90000000000001
$ws.Range("T5").Value = 0.689158453373768
$ws.Range("U5").Value = 242.5
$ws.Range("V5").Value = 0.07375752783015999
$ws.Range("W5").Value = 0.2939168961800982
$ws.Range("X5").Value = 0.06396583442419193
$ws.Range("Y5").Value = 0.2299510617559063
$ws.Range("Z5").Value = 2.218642000406587
$ws.Range("AA5").Value = 0.2454347283306264
$ws.Range("AB5").Value = 0.06064221083306002
$ws.Range("AC5").Value = 0.1847925174975664
$ws.Range("AD5").Value = 418.8
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 418.8
$ws.Range("AG5").Value = 176.3
$ws.Range("AH5").Value = 0.1129876436626558
$ws.Range("AI5").Value = 0.2142966791178427
$ws.Range("AJ5").Value = 0.05089344995814209
$ws.Range("AK5").Value = 0.1029910036219185
$ws.Range("AL5").Value = 17.4
$ws.Range("AM5").Value = 17.4
$ws.Range("AN5").Value = 0.6332022981554279
$ws.Range("AO5").Value = 38.55747126436782
$ws.Range("AP5").Value = 0.2665557907469006
$ws.Range("AQ5").Value = 38.55747126436782
$ws.Range("B6").Value = "Porto Seguro S.A. (BOVESPA:PSSA3)"
$ws.Range("D6").Value = 0.0354
$ws.Range("E6").Value = 0.109
$ws.Range("F6").Value = -0.0207
$ws.Range("G6").Value = 0.1928790048577766
$ws.Range("H6").Value = 0.1928790048577766
$ws.Range("I6").Value = 0.1672984584281021
$ws.Range("J6").Value = 0.1080186094922169
$ws.Range("K6").Value = 293.7
$ws.Range("L6").Value = 0.08246988459270491
$ws.Range("M6").Value = 129.82
$ws.Range("N6").Value = 0.04272643496577146
$ws.Range("O6").Value = 0.4420156622403814
$ws.Range("P6").Value = 122.7
$ws.Range("Q6").Value = 0.04038309636650869
$ws.Range("R6").Value = 0.4177732379979571
$ws.Range("S6").Value = 7.11999999999999
$ws.Range("T6").Value = 0.05484517023571091
$ws.Range("U6").Value = 41.5
$ws.Range("V6").Value = 0.01365850447604002
$ws.Range("W6").Value = 0.152191936988289
$ws.Range("X6").Value = 0.05997947110027935
$ws.Range("Y6").Value = 0.0922124658880096
$ws.Range("Z6").Value = 2.155098335854766
$ws.Range("AA6").Value = 0.2327907255580225
$ws.Range("AB6").Value = 0.05975473031074879
$ws.Range("AC6").Value = 0.1730359952472737
$ws.Range("AD6").Value = 23.1
$ws.Range("AF6").Value = 23.1
$ws.Range("AG6").Value = -18.4
$ws.Range("AH6").Value = 0.0075453209211171
$ws.Range("AI6").Value = 0.01441317776252574
$ws.Range("AJ6").Value = -0.006092715231788079
$ws.Range("AK6").Value = -0.01178580579041763
$ws.Range("AL6").Value = 15.1
$ws.Range("AM6").Value = 15.1
$ws.Range("AN6").Value = 0.03919904972000679
$ws.Range("AO6").Value = 39.45695364238411
$ws.Range("AP6").Value = -0.03122348549126082
$ws.Range("AQ6").Value = 39.45695364238411
$ws.Range("A7").Value = "Brazil"
$ws.Range("B7").Value = "Alper Consultoria e Corretora de Seguros S.A. (BOVESPA:APER3)"
$ws.Range("C7").Value = "Insurance (General)"
$ws.Range("D7").Value = -0.131
$ws.Range("G7").Value = 0.2106145251396648
$ws.Range("H7").Value = 0.2106145251396648
$ws.Range("I7").Value = 0.1435754189944134
$ws.Range("J7").Value = 0.1435754189944134
$ws.Range("K7").Value = 0.018
$ws.Range("L7").Value = 0.001005586592178771
$ws.Range("M7").Value = 0.1088
$ws.Range("N7").Value = 0.001528089887640449
$ws.Range("O7").Value = 6.044444444444443
$ws.Range("P7").Value = 0.1008
$ws.Range("Q7").Value = 0.001415730337078651
$ws.Range("R7").Value = 5.6
$ws.Range("S7").Value = 0.007999999999999993
$ws.Range("T7").Value = 0.07352941176470583
$ws.Range("U7").Value = 0.637
$ws.Range("V7").Value = 0.008946629213483145
$ws.Range("W7").Value = 0.0007725321888412016
$ws.Range("X7").Value = 0.06017144187425807
$ws.Range("Y7").Value = -0.05939890968541688
$ws.Range("Z7").Value = -12.2100954979536
$ws.Range("AA7").Value = -1.75306957708049
$ws.Range("AB7").Value = 0.05983338341570245
$ws.Range("AC7").Value = -1.812902960496192
$ws.Range("AD7").Value = 0.952
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0.952
$ws.Range("AG7").Value = 0.3149999999999999
$ws.Range("AH7").Value = 0.01319436744650183
$ws.Range("AI7").Value = 0.02897844880068185
$ws.Range("AJ7").Value = 0.004404670348877857
$ws.Range("AK7").Value = 0.009778053701691759
$ws.Range("AL7").Value = 0.07099999999999999
$ws.Range("AM7").Value = -0.234
$ws.Range("AN7").Value = 0.2659217877094972
$ws.Range("AO7").Value = 36.19718309859155
$ws.Range("AP7").Value = 0.08798882681564243
$ws.Range("AQ7").Value = -10.98290598290598
